# 23/04/2025 - cleanup test/demo data left over from switching the Excel
# automation backend away from win32com.client, plus a few formatting
# touch-ups (currency format for the "PREÇO" column, border fix on K2,
# wider "EMISSAO" column, updated zoom/selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 sample data clean-up -------------------------------------------
# A2: fix the capitalisation of the sample name
$ws.Range("A2").Value = "josé"
# B2: swap the placeholder CPF for a different dummy value
$ws.Range("B2").Value = "100.100.100-01"
# C2: used to hold a pasted pip-install console log by accident; clear it
$ws.Range("C2").Value = ""

# --- Formatting ------------------------------------------------------------
# K2 should pick up the same bordered style as the rest of the row (D2:I2)
$ws.Range("D2").Copy()
$ws.Range("K2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# PREÇO column (J) gets a real currency number format
$ws.Range("J2:J25").NumberFormat = """R$""#,##0.00"

# Widen the EMISSAO column (G) - closest achievable width to the 25.140625
# target (column widths are pixel-quantised by Excel's width model)
$ws.Columns("G").ColumnWidth = 25.0

# Update selection + zoom to match the saved view
[void]$ws.Range("C2").Select()
$excel.ActiveWindow.Zoom = 82
